$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.131385087966919
$ws.Range("B1").Value = 2.638291835784912
$ws.Range("C1").Value = 2.757874250411987
$ws.Range("D1").Value = 2.984084129333496
$ws.Range("E1").Value = 0.7646491527557373
